$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns for season record, matching the formatting of the
# existing header row (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill season record (Wins/Losses/Ties) for every player row
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
